# feat: add 2022-Q1 data
#
# 1. Insert a new sheet "2022-Q1" right before the "总计" (totals) sheet,
#    and populate it with the per-fund holding detail for 2022-Q1.
# 2. Insert a new summary row at the top of the "总计" sheet's data
#    (right after the header row) with the 2022-Q1 aggregate figures,
#    pushing the existing rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: add the "2022-Q1" sheet just before "总计"
# ---------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item("2021-Q4")

$ws = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$ws.Name = "2022-Q1"

# Copy header formatting (bold / border / centered) from an existing
# per-fund sheet so the new sheet matches the established look.
$srcSheet.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$srcSheet.Range("A2:A13").Copy()
$ws.Range("A2:A13").PasteSpecial(-4122)

# Header row
$ws.Range("B1").Value = "'基金代码"
$ws.Range("C1").Value = "'基金名称"
$ws.Range("D1").Value = "'基金规模"
$ws.Range("E1").Value = "'股票总仓位"
$ws.Range("F1").Value = "'仓位占比"
$ws.Range("G1").Value = "'持有市值(亿元)"
$ws.Range("H1").Value = "'仓位排名"

# Data rows (fund code / name are text; size, stock position, position
# ratio and held value are text-formatted numeric strings in the
# source data, so they are entered with a leading apostrophe to keep
# them as text instead of being auto-converted to numbers; rank is a
# real number).
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "'000727"
$ws.Range("C2").Value = "'融通健康产业灵活配置混合A"
$ws.Range("D2").Value = "'15.30"
$ws.Range("E2").Value = "'94.68"
$ws.Range("F2").Value = "'9.43"
$ws.Range("G2").Value = "'1.4428"
$ws.Range("H2").Value = 2

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "'009693"
$ws.Range("C3").Value = "'富国积极成长一年定期开放混合"
$ws.Range("D3").Value = "'17.82"
$ws.Range("E3").Value = "'98.74"
$ws.Range("F3").Value = "'5.48"
$ws.Range("G3").Value = "'0.9765"
$ws.Range("H3").Value = 1

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "'005984"
$ws.Range("C4").Value = "'兴业聚华混合A"
$ws.Range("D4").Value = "'21.22"
$ws.Range("E4").Value = "'24.30"
$ws.Range("F4").Value = "'1.82"
$ws.Range("G4").Value = "'0.3862"
$ws.Range("H4").Value = 6

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "'009274"
$ws.Range("C5").Value = "'融通健康产业灵活配置混合C"
$ws.Range("D5").Value = "'3.16"
$ws.Range("E5").Value = "'94.68"
$ws.Range("F5").Value = "'9.43"
$ws.Range("G5").Value = "'0.2980"
$ws.Range("H5").Value = 2

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "'000717"
$ws.Range("C6").Value = "'融通转型三动力灵活配置混合A"
$ws.Range("D6").Value = "'3.83"
$ws.Range("E6").Value = "'94.89"
$ws.Range("F6").Value = "'5.23"
$ws.Range("G6").Value = "'0.2003"
$ws.Range("H6").Value = 5

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "'161605"
$ws.Range("C7").Value = "'融通蓝筹成长混合"
$ws.Range("D7").Value = "'4.82"
$ws.Range("E7").Value = "'71.70"
$ws.Range("F7").Value = "'3.67"
$ws.Range("G7").Value = "'0.1769"
$ws.Range("H7").Value = 5

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "'005985"
$ws.Range("C8").Value = "'兴业聚华混合C"
$ws.Range("D8").Value = "'7.68"
$ws.Range("E8").Value = "'24.30"
$ws.Range("F8").Value = "'1.82"
$ws.Range("G8").Value = "'0.1398"
$ws.Range("H8").Value = 6

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "'005043"
$ws.Range("C9").Value = "'国寿安保健康科学混合A"
$ws.Range("D9").Value = "'0.99"
$ws.Range("E9").Value = "'85.72"
$ws.Range("F9").Value = "'4.32"
$ws.Range("G9").Value = "'0.0428"
$ws.Range("H9").Value = 4

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "'005044"
$ws.Range("C10").Value = "'国寿安保健康科学混合C"
$ws.Range("D10").Value = "'0.87"
$ws.Range("E10").Value = "'85.72"
$ws.Range("F10").Value = "'4.32"
$ws.Range("G10").Value = "'0.0376"
$ws.Range("H10").Value = 4

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "'009828"
$ws.Range("C11").Value = "'融通转型三动力灵活配置混合C"
$ws.Range("D11").Value = "'0.59"
$ws.Range("E11").Value = "'94.89"
$ws.Range("F11").Value = "'5.23"
$ws.Range("G11").Value = "'0.0309"
$ws.Range("H11").Value = 5

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "'009502"
$ws.Range("C12").Value = "'国寿安保创新医药股票A"
$ws.Range("D12").Value = "'0.54"
$ws.Range("E12").Value = "'81.60"
$ws.Range("F12").Value = "'4.20"
$ws.Range("G12").Value = "'0.0227"
$ws.Range("H12").Value = 2

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "'009503"
$ws.Range("C13").Value = "'国寿安保创新医药股票C"
$ws.Range("D13").Value = "'0.20"
$ws.Range("E13").Value = "'81.60"
$ws.Range("F13").Value = "'4.20"
$ws.Range("G13").Value = "'0.0084"
$ws.Range("H13").Value = 2

# ---------------------------------------------------------------------
# Step 2: insert the 2022-Q1 summary row at the top of "总计"'s data
# ---------------------------------------------------------------------
# Re-fetch the sheet by name: inserting "2022-Q1" above shifted "总计"'s
# position, and sheet object references are position-bound, so a
# reference captured before the insert would now resolve to the wrong
# sheet.
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# The insert copies row 1 (header) formatting onto the new row; reset
# it to plain/no style, then copy the real data-row formatting (from
# what is now row 3, the old row 2) onto the new row 2 instead.
$totalSheet.Range("A2:D2").Style = "Normal"
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 12
$totalSheet.Range("D2").Value = 3.76

# The "A" column is a plain sequential row index (0-based), not a
# formula, so re-number the rows that were pushed down one position.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
